$d = $word.ActiveDocument

# In the "Correções Banca" checklist every resolved item is struck
# through (the still-pending ones stay red/FF0000). The last item in
# the list - "falar mais sobre Qt e menos sobre alternativas (falar
# mais sobre Qt, tipo 'Como Qt trabalha com multiplataforma?' e
# emit)" - was still shown in plain/unformatted text. Mark it as
# resolved too, matching the rest of the list, by striking through the
# whole paragraph (including its paragraph mark, so the list bullet
# itself also carries the strike formatting like its siblings).

$needle = "mais sobre Qt e menos sobre alternativas"

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like ("*" + $needle + "*")) {
        $target = $p
    }
}

if ($target -eq $null) {
    # Fallback: the paragraph is the very last one in the document.
    $target = $d.Paragraphs.Last
}

$target.Range.Font.StrikeThrough = 1
